$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 94.5
$ws.Range("I4").Value = 126.28571
$ws.Range("K4").Value = 126.28571
$ws.Range("M4").Value = -12.28570999999999
$ws.Range("H9").Value = 6605.5
$ws.Range("I9").Value = 8690
$ws.Range("J9").Value = 3478.75
$ws.Range("K9").Value = 8690
$ws.Range("L9").Value = 3478.75
$ws.Range("M9").Value = -8521
$ws.Range("N9").Value = -3816.75
$ws.Range("H18").Value = 1038.4445
$ws.Range("I18").Value = 430.75
$ws.Range("K18").Value = 430.75
$ws.Range("M18").Value = -146.75
$ws.Range("H39").Value = 798.61536
$ws.Range("I39").Value = 34.727272
$ws.Range("K39").Value = 104.181816
$ws.Range("M39").Value = 191.818184
$ws.Range("H112").Value = 5141.6
$ws.Range("J112").Value = 5141.6
$ws.Range("L112").Value = 15424.8
$ws.Range("N112").Value = -17640.8
$ws.Range("H116").Value = 8971.951999999999
$ws.Range("I116").Value = 8366.25
$ws.Range("K116").Value = 8366.25
$ws.Range("M116").Value = -4924.25
$ws.Range("H138").Value = 4733.918
$ws.Range("I138").Value = 1381.72
$ws.Range("J138").Value = 7061.8335
$ws.Range("K138").Value = 4145.16
$ws.Range("L138").Value = 21185.5005
$ws.Range("M138").Value = 994.8400000000001
$ws.Range("N138").Value = -31465.5005

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2737.1875
$ws.Range("J2").Value = 2788
$ws.Range("L2").Value = 2788
$ws.Range("N2").Value = -3014
$ws.Range("H32").Value = 3638.875
$ws.Range("I32").Value = 3110.224
$ws.Range("K32").Value = 3110.224
$ws.Range("M32").Value = -2823.224
$ws.Range("H45").Value = 1431
$ws.Range("I45").Value = 842.4545000000001
$ws.Range("J45").Value = 3049.5
$ws.Range("K45").Value = 842.4545000000001
$ws.Range("L45").Value = 3049.5
$ws.Range("M45").Value = -465.4545000000001
$ws.Range("N45").Value = -3803.5
$ws.Range("H110").Value = 9372.916999999999
$ws.Range("I110").Value = 9348.223
$ws.Range("J110").Value = 9447
$ws.Range("K110").Value = 9348.223
$ws.Range("L110").Value = 9447
$ws.Range("M110").Value = -7303.223
$ws.Range("N110").Value = -13537
$ws.Range("H116").Value = 2737.1875
$ws.Range("J116").Value = 2788
$ws.Range("L116").Value = 2788
$ws.Range("N116").Value = -7376
$ws.Range("H132").Value = 100269770
$ws.Range("I132").Value = 55734.6
$ws.Range("J132").Value = 200483800
$ws.Range("K132").Value = 167203.8
$ws.Range("L132").Value = 601451400
$ws.Range("M132").Value = -164673.8
$ws.Range("N132").Value = -601456460

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2737.1875
$ws.Range("J3").Value = 2788
$ws.Range("L3").Value = 2788
$ws.Range("N3").Value = -3016
$ws.Range("H82").Value = 10751.571
$ws.Range("J82").Value = 47783
$ws.Range("L82").Value = 47783
$ws.Range("N82").Value = -48549
$ws.Range("H85").Value = 10751.571
$ws.Range("J85").Value = 47783
$ws.Range("L85").Value = 47783
$ws.Range("N85").Value = -50435
$ws.Range("H86").Value = 21830.334
$ws.Range("I86").Value = 36161.668
$ws.Range("K86").Value = 36161.668
$ws.Range("M86").Value = -35038.668
$ws.Range("H89").Value = 21830.334
$ws.Range("I89").Value = 36161.668
$ws.Range("K89").Value = 180808.34
$ws.Range("M89").Value = -175192.34
$ws.Range("H94").Value = 1708.6
$ws.Range("I94").Value = 1506.6154
$ws.Range("J94").Value = 2083.7144
$ws.Range("K94").Value = 1506.6154
$ws.Range("L94").Value = 2083.7144
$ws.Range("M94").Value = -1055.6154
$ws.Range("N94").Value = -2985.7144
$ws.Range("H105").Value = 6181.189
$ws.Range("I105").Value = 8580.4375
$ws.Range("K105").Value = 8580.4375
$ws.Range("M105").Value = -6833.4375

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 8423.071
$ws.Range("I22").Value = 11407.777
$ws.Range("K22").Value = 11407.777
$ws.Range("M22").Value = -11057.777
$ws.Range("H99").Value = 7946.6924
$ws.Range("I99").Value = 8098.222
$ws.Range("J99").Value = 7605.75
$ws.Range("K99").Value = 8098.222
$ws.Range("L99").Value = 7605.75
$ws.Range("M99").Value = -6600.222
$ws.Range("N99").Value = -10601.75
$ws.Range("H126").Value = 7946.6924
$ws.Range("I126").Value = 8098.222
$ws.Range("J126").Value = 7605.75
$ws.Range("K126").Value = 24294.666
$ws.Range("L126").Value = 22817.25
$ws.Range("M126").Value = -21824.666
$ws.Range("N126").Value = -27757.25
$ws.Range("H131").Value = 24066.666
$ws.Range("J131").Value = 28600
$ws.Range("L131").Value = 28600
$ws.Range("N131").Value = -38680
$ws.Range("H132").Value = 203853.6
$ws.Range("I132").Value = 288787
$ws.Range("K132").Value = 866361
$ws.Range("M132").Value = -863831
$ws.Range("H141").Value = 141598
$ws.Range("J141").Value = 169497.5
$ws.Range("L141").Value = 169497.5
$ws.Range("N141").Value = -179857.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 955.55554
$ws.Range("J23").Value = 1228.3334
$ws.Range("L23").Value = 3685.0002
$ws.Range("N23").Value = -4155.0002
$ws.Range("H131").Value = 1156.6428
$ws.Range("J131").Value = 1757.4166
$ws.Range("L131").Value = 5272.2498
$ws.Range("N131").Value = -15352.2498
$ws.Range("H133").Value = 4080.625
$ws.Range("I133").Value = 1315.8
$ws.Range("J133").Value = 8688.666999999999
$ws.Range("K133").Value = 3947.4
$ws.Range("L133").Value = 26066.001
$ws.Range("M133").Value = 1112.6
$ws.Range("N133").Value = -36186.001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3968.3333
$ws.Range("I113").Value = 3933.6667
$ws.Range("J113").Value = 4037.6667
$ws.Range("K113").Value = 3933.6667
$ws.Range("L113").Value = 4037.6667
$ws.Range("M113").Value = -1763.6667
$ws.Range("N113").Value = -8377.6667
$ws.Range("H122").Value = 41668388
$ws.Range("I122").Value = 1463.7
$ws.Range("K122").Value = 4391.1
$ws.Range("M122").Value = -1941.1
$ws.Range("H126").Value = 2544.7144
$ws.Range("I126").Value = 2468.8333
$ws.Range("K126").Value = 7406.499899999999
$ws.Range("M126").Value = -4936.499899999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H123").Value = 68662.836
$ws.Range("J123").Value = 74995.39999999999
$ws.Range("L123").Value = 74995.39999999999
$ws.Range("N123").Value = -84795.39999999999
$ws.Range("H131").Value = 74644.2
$ws.Range("J131").Value = 89308.336
$ws.Range("L131").Value = 89308.336
$ws.Range("N131").Value = -99388.336
$ws.Range("H132").Value = 4682.027
$ws.Range("I132").Value = 2235.48
$ws.Range("J132").Value = 9779
$ws.Range("K132").Value = 6706.440000000001
$ws.Range("L132").Value = 29337
$ws.Range("M132").Value = -4176.440000000001
$ws.Range("N132").Value = -34397

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1738.6842
$ws.Range("I107").Value = 648.8889
$ws.Range("J107").Value = 2719.5
$ws.Range("K107").Value = 1946.6667
$ws.Range("L107").Value = 8158.5
$ws.Range("M107").Value = -26.66670000000022
$ws.Range("N107").Value = -11998.5
$ws.Range("H113").Value = 1282.125
$ws.Range("I113").Value = 1449.6
$ws.Range("J113").Value = 1003
$ws.Range("K113").Value = 4348.799999999999
$ws.Range("L113").Value = 3009
$ws.Range("M113").Value = -2178.799999999999
$ws.Range("N113").Value = -7349
$ws.Range("H129").Value = 40000
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()
$ws.Range("H131").Value = 87999
$ws.Range("J131").Value = 87999
$ws.Range("L131").Value = 87999
$ws.Range("N131").Value = -98079
